$d = $word.ActiveDocument

$marker = "This is a totally new comment."

# Locate the paragraph that currently holds the marker text (and, along with
# it, the _GoBack bookmark Word leaves at the last edit position).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq ($marker + [char]13)) {
        $target = $p
    }
}
if ($target -eq $null) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like ($marker + "*")) {
            $target = $p
        }
    }
}

# Temporarily append a one-character placeholder right after the existing
# text (i.e. before the bookmark). This keeps the upcoming paragraph split
# away from the exact boundary position the bookmark already occupies, so
# the split lands cleanly and the bookmark ends up carried into the new
# paragraph instead of staying behind with the old text.
$target.Range.InsertAfter("X") | Out-Null

$splitPos = $target.Range.Start + $marker.Length
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter() | Out-Null

# The freshly created paragraph now contains just "X" followed by the
# relocated _GoBack bookmark.
$newPara = $target.Next()
$fullNewRange = $newPara.Range

$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
       "<w:r><w:t xml:space='preserve'>This is yet another comment. This comment should not make into the repo b/c </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/>" +
       "<w:r><w:t>lmiksa-gorillagroup</w:t></w:r>" +
       "<w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> is no longer a collaborator.</w:t></w:r>" +
       "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
       "<w:bookmarkEnd w:id='0'/>" +
       "</w:p>"

# Replace the whole "X" + bookmark paragraph content with the real comment
# text/runs/proofErr markers, re-emitting the _GoBack bookmark at the end.
$fullNewRange.InsertXML($xml) | Out-Null
